# 🚌 141: 31/12 12:00 LP1912+6203+6173
# Appends newly-scraped schedule rows to the three worksheets and refreshes
# the "Última actualización" / "Total filas" header cells on each sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "LP1912"  -> columns: A(-), B=Hora_Scrap, C=Hora_Llegada,
#                        D=Línea, E=Minutos, F=Parada, G=Fecha
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 31/12/2025 09:00:57"
$ws1.Range("A3").Value = "Total filas: 810"

$sheet1Rows = @(
    @("", "09:00:46", "09:03", "23_HERNANDEZ", 3, "LP1912", "31/12/2025"),
    @("", "09:00:46", "09:08", "16_SANTA ANA", 8, "LP1912", "31/12/2025"),
    @("", "09:00:46", "09:14", "11_ETCHEVERRY", 14, "LP1912", "31/12/2025"),
    @("", "09:00:46", "09:16", "27_EL RETIRO", 16, "LP1912", "31/12/2025"),
    @("", "09:00:46", "09:21", "16_SANTA ANA", 21, "LP1912", "31/12/2025"),
    @("", "09:00:46", "09:26", "215_EL PELIGRO", 26, "LP1912", "31/12/2025"),
    @("", "09:00:46", "09:27", "23_HERNANDEZ", 27, "LP1912", "31/12/2025"),
    @("", "09:00:46", "09:33", "23_HERNANDEZ", 33, "LP1912", "31/12/2025"),
    @("", "09:00:46", "09:41", "16_SANTA ANA", 41, "LP1912", "31/12/2025"),
    @("", "09:00:46", "09:44", "14_ABASTO", 44, "LP1912", "31/12/2025"),
    @("", "09:00:46", "09:51", "15_ABASTO", 51, "LP1912", "31/12/2025"),
    @("", "09:00:46", "09:54", "10_OLMOS", 54, "LP1912", "31/12/2025"),
    @("", "09:00:46", "10:02", "215C_EL PATO", 62, "LP1912", "31/12/2025"),
    @("", "09:00:46", "10:04", "14_ABASTO", 64, "LP1912", "31/12/2025"),
    @("", "09:00:46", "10:14", "10_OLMOS", 74, "LP1912", "31/12/2025"),
    @("", "09:00:46", "10:24", "11_ETCHEVERRY", 84, "LP1912", "31/12/2025"),
    @("", "09:00:46", "10:26", "15X38_ABASTO", 86, "LP1912", "31/12/2025"),
    @("", "09:00:46", "10:34", "10_OLMOS", 94, "LP1912", "31/12/2025"),
    @("", "09:00:46", "10:37", "16_P MOR-SANTA ANA", 97, "LP1912", "31/12/2025")
)

$startRow1 = 793
for ($i = 0; $i -lt $sheet1Rows.Count; $i++) {
    $rowData = $sheet1Rows[$i]
    $r = $startRow1 + $i
    for ($j = 0; $j -lt $rowData.Count; $j++) {
        $ws1.Cells.Item($r, $j + 1).Value = $rowData[$j]
    }
}

# ---------------------------------------------------------------------
# Sheet 2: "LP1912-215" -> columns: A(-), B=Fecha, C=Hora_Scrap,
#                           D=Hora_Llegada, E=Línea, F=Minutos, G=Parada
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 31/12/2025 09:00:57"
$ws2.Range("A3").Value = "Total filas: 60"

$sheet2Rows = @(
    @("", "31/12/2025", "09:00:46", "09:26", "215_EL PELIGRO", 26, "LP1912"),
    @("", "31/12/2025", "09:00:46", "10:02", "215C_EL PATO", 62, "LP1912")
)

$startRow2 = 60
for ($i = 0; $i -lt $sheet2Rows.Count; $i++) {
    $rowData = $sheet2Rows[$i]
    $r = $startRow2 + $i
    for ($j = 0; $j -lt $rowData.Count; $j++) {
        $ws2.Cells.Item($r, $j + 1).Value = $rowData[$j]
    }
}

# ---------------------------------------------------------------------
# Sheet 3: "6203-6173" -> columns: A(-), B=Fecha, C=Hora_Scrap,
#                          D=Hora_Llegada, E=Línea, F=Minutos, G=Parada
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 31/12/2025 09:00:57"
$ws3.Range("A3").Value = "Total filas: 96"

$sheet3Rows = @(
    @("", "31/12/2025", "09:00:51", "10:08", "215C_LA PLATA", 68, "L6203"),
    @("", "31/12/2025", "09:00:56", "10:22", "215A_LA PLATA", 82, "L6173"),
    @("", "31/12/2025", "09:00:56", "10:30", "215B_LP-P MOR-1 Y 57", 90, "L6173")
)

$startRow3 = 95
for ($i = 0; $i -lt $sheet3Rows.Count; $i++) {
    $rowData = $sheet3Rows[$i]
    $r = $startRow3 + $i
    for ($j = 0; $j -lt $rowData.Count; $j++) {
        $ws3.Cells.Item($r, $j + 1).Value = $rowData[$j]
    }
}
